# Fixed several issues associated with quantifying r in the DDE model for
# species with highly nonlinear trait responses.
#
# - Corrects the rMax (E) value for two Clavigralla tomentosicollis
#   populations (Nigeria, Burkina Faso) whose trait responses were highly
#   nonlinear, which also recalculates all of the dependent ratio/delta
#   formulas in columns J:Q for those rows.
# - Highlights (yellow fill) the raw trait-response columns (r.TPC.h,
#   r.TPC.f, r.model.h, r.model.f) for the affected / flagged rows so they
#   are easy to spot.
# - Updates the active selection left on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct rMax values for the two affected species ---------------------
$ws.Range("E4").Value = 0.088
$ws.Range("E5").Value = 0.088

# --- Highlight the highly-nonlinear trait-response rows --------------------
$yellow = 65535   # RGB(255,255,0) packed as BGR for the COM Color property

$ws.Range("F4:I5").Interior.Color = $yellow
$ws.Range("F17:I17").Interior.Color = $yellow
$ws.Range("F23:I23").Interior.Color = $yellow

# --- Leave the selection where the author left it --------------------------
$ws.Range("E6").Select()
